$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Update the publish-date label in H9
$ws.Range("H9").Value = "1402-01-28 (8)"

# Rows where D:H are all reset to 0
$zeroRows = @(11,12,13,14,16,17,18,19,20,21,22,24,27)
foreach ($r in $zeroRows) {
    foreach ($col in @("D","E","F","G","H")) {
        $ws.Range("$col$r").Value = 0
    }
}

# Rows where D:H are all reset to "-"
$dashRows = @(15,23)
foreach ($r in $dashRows) {
    foreach ($col in @("D","E","F","G","H")) {
        $ws.Range("$col$r").Value = "-"
    }
}

# Row 25: D=0, E="-", F="-", G=0, H=0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = "-"
$ws.Range("F25").Value = "-"
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0

# Row 26: D=0, E="-", F=0, G=0, H=0
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = "-"
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 0
